$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --------------------------------------------------------------------
# 1) Preserve the existing "source" block (currently at rows 26-27) by
#    writing its text further down to rows 32-33 BEFORE it gets
#    overwritten by the new MSME definition table in step 2.
# --------------------------------------------------------------------
$sourceTitle = $ws.Range("A26").Value2
$sourceText  = $ws.Range("A27").Value2

$ws.Range("A32").Value = $sourceTitle
$ws.Range("A32").Font.Bold = $true

$ws.Range("A33").Value = $sourceText
$ws.Range("A33").Font.Italic = $true

# --------------------------------------------------------------------
# 2) Insert the new "MSME definition" table in rows 23-27.
# --------------------------------------------------------------------

# Header row (bold), like the other table headers on the sheet
$ws.Range("B23").Value = "Number of employees"
$ws.Range("C23").Value = "Assets (local currency, unless noted otherwise)"
$ws.Range("D23").Value = "Turnover (local currency, unless noted otherwise)"
$ws.Range("B23:D23").Font.Bold = $true

# Micro
$ws.Range("A24").Value = "Micro"
$ws.Range("B24").Value = "<10"
$ws.Range("C24").Value = "'"
$ws.Range("D24").Value = "'"

# Small
$ws.Range("A25").Value = "Small"
$ws.Range("B25").Value = "<50"
$ws.Range("C25").Value = "'"
$ws.Range("D25").Value = "'"

# Medium (A26 previously held the bold "SME Performance Review EU" text,
# so clear the bold formatting now that it holds plain data)
$ws.Range("A26").Value = "Medium"
$ws.Range("A26").Font.Bold = $false
$ws.Range("B26").Value = "<250"
$ws.Range("C26").Value = "'"
$ws.Range("D26").Value = "'"

# Large (A27 previously held the italic source text, so clear the italics)
$ws.Range("A27").Value = "Large"
$ws.Range("A27").Font.Italic = $false
$ws.Range("B27").Value = ">249"
$ws.Range("C27").Value = "'"
$ws.Range("D27").Value = "'"
